# "Generate Report for Handback" -- mark the two in-flight localization
# files as handed back (in sync with en-US) and record the Latest Target
# File / Latest Handback File / Latest Handback DateTime for each locale.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Source-file link targets (unchanged -- reused for the new "Latest Target
# File" column, which simply repeats the source-file hyperlink).
$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/048c823da19957da67982ea9af74370f81b5ea46/e2e/e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/048c823da19957da67982ea9af74370f81b5ea46/e2e/e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/048c823da19957da67982ea9af74370f81b5ea46/.localization-config"

$mdName1 = "e369dccd-7b85-4c68-a7da-f277bab9fa44.md"
$mdName2 = "e52524ed-a6eb-4c40-9710-5e56e59168b0.md"
$cfgName = ".localization-config"

# ------------------------------------------------------------------
# Overview sheet: just the status text changes (same text used for both
# the zh-cn and de-de columns of each source-file row).
# ------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("B3").Value = $statusHandedBack
$ov.Range("C3").Value = $statusHandedBack

# ------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de share the same layout
#   A Source File Name   B Status            C Latest Handoff File
#   D Latest Handoff Datetime   E Latest Target File   F Latest Handback File
#   G Latest Handback DateTime  H Handoff Reason        I Dependency From
# ------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn";
       XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e9511738520c68f20f8bd717f513003b7673411/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e9511738520c68f20f8bd717f513003b7673411/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf";
       XlfName1 = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.zh-cn.xlf";
       XlfName2 = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.zh-cn.xlf";
       Handback = "2016-03-09 09:44:25" },
    @{ Name = "de-de";
       XlfUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2985ce4d130740c6fdcd978887d06b7f6e188f16/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf";
       XlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2985ce4d130740c6fdcd978887d06b7f6e188f16/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf";
       XlfName1 = "e369dccd-7b85-4c68-a7da-f277bab9fa44.f79e141b90cc523ca1fb057b3f77fc7bc7b90aee.de-de.xlf";
       XlfName2 = "e52524ed-a6eb-4c40-9710-5e56e59168b0.8c3861640126a5d8d85b66e2fe82398eb991207a.de-de.xlf";
       Handback = "2016-03-09 09:44:41" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Name)

    # Status column: both handed-off rows are now synced with en-US.
    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("B3").Value = $statusHandedBack

    # Latest Handback DateTime for the two rows that were just handed back.
    $ws.Range("G2").Value = $loc.Handback
    $ws.Range("G3").Value = $loc.Handback

    # Rebuild the hyperlinks in final left-to-right / top-to-bottom order so
    # the newly added Latest Target File (E) and Latest Handback File (F)
    # links land between the existing columns, exactly where they belong.
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("C2"), $loc.XlfUrl1, "", "", $loc.XlfName1)
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl1, "", "", $mdName1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $loc.XlfUrl1, "", "", $loc.XlfName1)

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("C3"), $loc.XlfUrl2, "", "", $loc.XlfName2)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl2, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $loc.XlfUrl2, "", "", $loc.XlfName2)

    $ws.Hyperlinks.Add($ws.Range("A4"), $cfgUrl, "", "", $cfgName)

    # Give every hyperlinked cell the workbook's existing hyperlink look
    # (underline + custom blue) instead of Excel's default theme color.
    foreach ($addr in @("A2", "C2", "E2", "F2", "A3", "C3", "E3", "F3", "A4")) {
        $rng = $ws.Range($addr)
        $rng.Font.Underline = 2
        $rng.Font.Color = 15570276
    }
}
